$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.693.54"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "3.533.19"
$ws.Range("E3").Value = "  +6.14%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'238.82"
$ws.Range("E5").Value = "  +3.69%  "
$ws.Range("D6").Value = "'632.99"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("D7").Value = "'1.44"
$ws.Range("E7").Value = "  +6.46%  "
$ws.Range("D8").Value = "'0.398"
$ws.Range("E8").Value = "  +4.65%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +9.16%  "
$ws.Range("D11").Value = "3.525.58"
$ws.Range("E11").Value = "  +5.86%  "
$ws.Range("D12").Value = "'43.61"
$ws.Range("E12").Value = "  +4.86%  "
$ws.Range("E13").Value = "  +5.48%  "
$ws.Range("D14").Value = "'6.28"
$ws.Range("E14").Value = "  +5.44%  "
$ws.Range("D15").Value = "4.202.63"
$ws.Range("E15").Value = "  +6.22%  "
$ws.Range("D16").Value = "94.580.28"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("E17").Value = "  +4.43%  "
$ws.Range("D18").Value = "'8.31"
$ws.Range("E18").Value = "  +5.10%  "
$ws.Range("D19").Value = "3.532.32"
$ws.Range("E19").Value = "  +6.05%  "
$ws.Range("D20").Value = "'13.06"
$ws.Range("E20").Value = "  +19.50%  "
$ws.Range("D21").Value = "'18.07"
$ws.Range("E21").Value = "  +5.84%  "
$ws.Range("D22").Value = "'0.497"
$ws.Range("E22").Value = "  +11.16%  "
$ws.Range("D23").Value = "'515.69"
$ws.Range("E23").Value = "  +5.68%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "'6.74"
$ws.Range("E25").Value = "  +10.02%  "
$ws.Range("D26").Value = "'0.0000189"
$ws.Range("E26").Value = "  +5.57%  "
$ws.Range("D27").Value = "'92.74"
$ws.Range("E27").Value = "  +3.80%  "
$ws.Range("D28").Value = "'12.24"
$ws.Range("E28").Value = "  +6.50%  "
$ws.Range("D29").Value = "'3.06"
$ws.Range("E29").Value = "  +16.11%  "
$ws.Range("D30").Value = "'11.54"
$ws.Range("E30").Value = "  +5.21%  "
$ws.Range("E31").Value = "  +6.05%  "
$ws.Range("D33").Value = "'0.183"
$ws.Range("E33").Value = "  +7.25%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D35").Value = "'30.03"
$ws.Range("E35").Value = "  +6.66%  "
$ws.Range("E36").Value = "  +7.42%  "
$ws.Range("D37").Value = "'583.97"
$ws.Range("E37").Value = "  +10.22%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'7.61"
$ws.Range("E38").Value = "  +4.68%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'1.45"
$ws.Range("E39").Value = "  +7.45%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'0.929"
$ws.Range("E41").Value = "  +7.35%  "
$ws.Range("D42").Value = "'0.151"
$ws.Range("E42").Value = "  +3.48%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0426"
$ws.Range("E43").Value = "  +5.65%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'23.77"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").Value = "'1.71"
$ws.Range("E45").Value = "  +3.79%  "
$ws.Range("D46").Value = "'5.56"
$ws.Range("E46").Value = "  +4.82%  "
$ws.Range("D47").Value = "'3.54"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'2.18"
$ws.Range("E48").Value = "  +4.44%  "
$ws.Range("D49").Value = "'54.02"
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("D50").Value = "'8.16"
$ws.Range("E50").Value = "  +4.54%  "
$ws.Range("D51").Value = "'3.10"
$ws.Range("E51").Value = "  +3.45%  "
